{"js": "// Prologue Text Script edits \u2014 Level 0.\n// Each change below narrows a stock phrase to the exact sentence that\n// changed in the source doc and rewrites it with the revised wording\n// (Word itself merges same-formatted adjacent runs back together, so a\n// plain search + insertText(\"Replace\") reproduces the same end state the\n// diff's run-splitting produced).\n\nconst body = context.document.body;\n\nasync function replaceOnce(findText, newText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"*You duck right before...\" -> \"*You ducked right before...\"\nawait replaceOnce(\n  \"*You duck right before the female government agent fires her handgun!\",\n  \"*You ducked right before the female government agent fires her handgun!\"\n);\n\n// 2) \"You: W-What? \" -> \"Player: W-What? \"\nawait replaceOnce(\n  \"You: W-What? \",\n  \"Player: W-What? \"\n);\n\n// 3) \"Their clearly concussed or something. \" -> \"They\u2019re clearly concussed or something. \"\nawait replaceOnce(\n  \"Their clearly concussed or something. \",\n  \"They\\u2019re clearly concussed or something. \"\n);\n\n// 4) \"...a man who has the sun behind him preventing you from seeing him clearly.\"\n//    -> \"...a feminine figure who has the sun behind them preventing you from seeing clearly.\"\nawait replaceOnce(\n  \"*You are forced to knell before a man who has the sun behind him preventing you from seeing him clearly.\",\n  \"*You are forced to knell before a feminine figure who has the sun behind them preventing you from seeing clearly.\"\n);\n\n// 5) \"...see the man reeling backwards manically laughing with his hand covering his face. \"\n//    -> \"...see the unknown female? reeling backwards manically laughing with their hand covering their face. \"\nawait replaceOnce(\n  \"*You glance up to see the man reeling backwards manically laughing with his hand covering his face. \",\n  \"*You glance up to see the unknown female? reeling backwards manically laughing with their hand covering their face. \"\n);\n", "ps1": "# Prologue Text Script edits \u2014 Level 0.\n# Each block below finds the exact sentence that changed in the source\n# document and replaces it with the revised wording. Word merges\n# identically-formatted adjacent runs back together on save, so a plain\n# Find/Replace over $d.Content reproduces the same end state that the\n# diff's run-splitting represents.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\n# 1) \"*You duck right before...\" -> \"*You ducked right before...\"\nReplace-Once \"*You duck right before the female government agent fires her handgun!\" \"*You ducked right before the female government agent fires her handgun!\"\n\n# 2) \"You: W-What? \" -> \"Player: W-What? \"\nReplace-Once \"You: W-What? \" \"Player: W-What? \"\n\n# 3) \"Their clearly concussed or something. \" -> \"They\u2019re clearly concussed or something. \"\nReplace-Once \"Their clearly concussed or something. \" \"They\u2019re clearly concussed or something. \"\n\n# 4) \"...a man who has the sun behind him preventing you from seeing him clearly.\"\n#    -> \"...a feminine figure who has the sun behind them preventing you from seeing clearly.\"\nReplace-Once \"*You are forced to knell before a man who has the sun behind him preventing you from seeing him clearly.\" \"*You are forced to knell before a feminine figure who has the sun behind them preventing you from seeing clearly.\"\n\n# 5) \"...see the man reeling backwards manically laughing with his hand covering his face. \"\n#    -> \"...see the unknown female? reeling backwards manically laughing with their hand covering their face. \"\nReplace-Once \"*You glance up to see the man reeling backwards manically laughing with his hand covering his face. \" \"*You glance up to see the unknown female? reeling backwards manically laughing with their hand covering their face. \"\n"}
